$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed set of test data (email / password / username triples).
$rows = @(
  @("cetsbvek@vmyxalzn.com",      "915LK2Ur",      "oqtxycgj"),
  @("eotwmrjhv@tiaccgqee.com",    "kOaJvHBBS",     "ddotrshah"),
  @("addssx@jeftdt.com",          "7pBVRm",        "jvbedq"),
  @("qtlmljnyzvt@ulafxwyahww.com","KGMDsglWWc1",   "xncirpsbvby"),
  @("rdmevr@ingsjg.com",          "L8gAcb",        "rgaufj"),
  @("rdm12evr@ingsjg.com",        "L8gA12cb",      "rga12ufj"),
  @("azinhdz@cmjeklf.com",        "bJU25Wa",       "dgglykg"),
  @("kwsxpiea@vrngtydk.com",      "9OgU4nxL",      "qxuecfcr"),
  @("wxkklqa@rnzbgle.com",        "p5K5x1Q",       "vgmmybe"),
  @("hncrcushlut@cfqqmicsczn.com","5QNoRcPQK4R",   "oeenujvywzu"),
  @("xwljgviytn@thizexnaki.com",  "6dBVBqp3kg",    "zujlpnxecw")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 1
  $ws.Cells.Item($r, 1).Value = $rows[$i][0]
  $ws.Cells.Item($r, 2).Value = $rows[$i][1]
  $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

# Row 6 (the injected/regenerated record) gets a hyperlink on its email cell.
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:rdm12evr@ingsjg.com")

$ws.Range("C6").Select() | Out-Null
